$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix OP-07 target value text: "500>" -> ">500"
$ws.Range("B13").Value = ">500"

# 2. Add new row 17 for OP-11 (ventas por canales digitales)

# A17: ID indicador cell -> border on all sides, left/center aligned, wrap text
$ws.Range("A7").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("A17").HorizontalAlignment = -4131
$ws.Range("A17").Value = "OP-11"

# B17: Valor objetivo -> border all sides, vertical centered, wrap text, general number format
$ws.Range("A7").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("B17").HorizontalAlignment = 1
$ws.Range("B17").Value = 100000

# C17: Valor inicial -> same style as B17
$ws.Range("A7").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("C17").HorizontalAlignment = 1
$ws.Range("C17").Value = 0

# D17: En lenguaje natural -> border on left/right only, wrap text
$ws.Range("D17").Borders.Item(7).LineStyle = 1
$ws.Range("D17").Borders.Item(7).Weight = 2
$ws.Range("D17").Borders.Item(10).LineStyle = 1
$ws.Range("D17").Borders.Item(10).Weight = 2
$ws.Range("D17").WrapText = $true
$ws.Range("D17").Value = "Queremos vender 10000 artículos mediante canales digitales en el último año"

$ws.Rows.Item(17).RowHeight = 30

# 3. Highlight D8 in yellow (same fill used elsewhere in the sheet)
$ws.Range("D8").Interior.Color = 65535

# 4. Widen column D slightly to fit the new text
$ws.Columns.Item(4).ColumnWidth = 44

# 5. Move the active selection to D8
$ws.Range("D8").Select()
